# Fixed naive component forecaster bug - Presentation state 11.02.
# The error table is recomputed: each existing row's B:G values shift down
# by one row (row N's values move to row N+1), the oldest row (row 11)
# drops off, and a brand new row of statistics is inserted at row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current B2:G10 values (these will be moved down to B3:G11)
$shiftRange = $ws.Range("B2:G10").Value2

# Write them into B3:G11
$ws.Range("B3:G11").Value2 = $shiftRange

# Write the new first row of statistics into B2:G2
$ws.Range("B2").Value2 = 0.002099636470939166
$ws.Range("C2").Value2 = 0.1232588647135871
$ws.Range("D2").Value2 = 0.03139795277152564
$ws.Range("E2").Value2 = 0.1771946747832046
$ws.Range("F2").Value2 = 0.1834010368550752
$ws.Range("G2").Value2 = 15
